$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览
$wsExhibit.Range("F4").Value = 5187
$wsExhibit.Range("F5").Value = 5187
$wsExhibit.Range("F7").Value = 162
$wsExhibit.Range("F8").Value = 213
$wsExhibit.Range("F11").Value = 177
$wsExhibit.Range("F12").Value = 8658
$wsExhibit.Range("F13").Value = 8658
$wsExhibit.Range("F14").Value = 276
$wsExhibit.Range("F15").Value = 138
$wsExhibit.Range("F16").Value = 631
$wsExhibit.Range("F18").Value = 2577
$wsExhibit.Range("F19").Value = 6333
$wsExhibit.Range("F20").Value = 2327
$wsExhibit.Range("G20").Value = 19.9
$wsExhibit.Range("F23").Value = 2538
$wsExhibit.Range("F25").Value = 19
$wsExhibit.Range("F26").Value = 6513
$wsExhibit.Range("F27").Value = 207
$wsExhibit.Range("F28").Value = 78
$wsExhibit.Range("F29").Value = 145
$wsExhibit.Range("F32").Value = 7061
$wsExhibit.Range("F40").Value = 1
$wsExhibit.Range("F41").Value = 41
$wsExhibit.Range("F43").Value = 2543
$wsExhibit.Range("F48").Value = 542
$wsExhibit.Range("F49").Value = 3196
$wsExhibit.Range("F51").Value = 1128

# 演出
$wsShow.Range("F7").Value = 84
$wsShow.Range("F10").Value = 14

# 全部类型
$wsAll.Range("F3").Value = 5187
$wsAll.Range("F4").Value = 5187
$wsAll.Range("F6").Value = 162
$wsAll.Range("F7").Value = 213
$wsAll.Range("F10").Value = 177
$wsAll.Range("F11").Value = 8658
$wsAll.Range("F12").Value = 8658
$wsAll.Range("F13").Value = 276
$wsAll.Range("F14").Value = 138
$wsAll.Range("F15").Value = 631
$wsAll.Range("F16").Value = 2577
$wsAll.Range("F19").Value = 6333
$wsAll.Range("F20").Value = 2327
$wsAll.Range("G20").Value = 19.9
$wsAll.Range("F21").Value = 84
$wsAll.Range("F22").Value = 2538
$wsAll.Range("F25").Value = 19
$wsAll.Range("F26").Value = 6513
$wsAll.Range("F27").Value = 207
$wsAll.Range("F28").Value = 14
$wsAll.Range("F29").Value = 78
$wsAll.Range("F30").Value = 145
$wsAll.Range("F33").Value = 7061
$wsAll.Range("F38").Value = 41
$wsAll.Range("F41").Value = 2543
$wsAll.Range("F45").Value = 542
$wsAll.Range("F47").Value = 3197
$wsAll.Range("F50").Value = 1128

Write-Output "Applied all changes"
